# Add a new "time" worksheet in front of the existing sheets and populate
# it with the time/frequency conversion-factor table, mirroring the layout
# used by the other conversion-factor sheets in this workbook.

$wb = $excel.ActiveWorkbook

# --- Create the new sheet as the very first tab -----------------------
$ws = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$ws.Name = "time"

# --- Header row ----------------------------------------------------------
$ws.Cells.Item(1,1).Value = "Source"
$ws.Cells.Item(1,2).Value = "Destination"
$ws.Cells.Item(1,3).Value = "Factor"
$ws.Cells.Item(1,4).Value = "Source name"
$ws.Cells.Item(1,5).Value = "Destination Name"

# --- Data rows -------------------------------------------------------------
# row, A(source), B(destination), C(factor), D(source name), E(destination name)
$rows = @(
    @(2,  "min", "s",   60,                 "second", "minute"),
    @(3,  "hour","min", 60,                 "hour",   "minute"),
    @(4,  "s",   "ms",  1000,               "second", "milisecond"),
    @(5,  "s",   "μs",  1000000,            "second", "microsecond"),
    @(6,  "s",   "ns",  1000000000,         "second", "nanosecond"),
    @(7,  "s",   "ps",  1000000000000,      "second", "picosecond"),
    @(8,  "s",   "fs",  1000000000000000,   "second", "femtpsecond"),
    @(9,  "Hz",  "kHz", (1/1000),           "Hertz",  "kilohertz"),
    @(10, "Hz",  "MHz", (1/1000000),        "Hertz",  "metahertz"),
    @(11, "Hz",  "GHz", (1/1000000000),     "Hertz",  "gigahertz"),
    @(12, "Hz",  "THz", (1/1000000000000),  "Hertz",  "terarhertz")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum,1).Value = $r[1]
    $ws.Cells.Item($rowNum,2).Value = $r[2]
    $ws.Cells.Item($rowNum,3).Value = $r[3]
    $ws.Cells.Item($rowNum,4).Value = $r[4]
    $ws.Cells.Item($rowNum,5).Value = $r[5]
}

# Scientific-notation number format for the Factor column (rows 2-12), plus
# a stray formatted-but-empty cell at C13 left over from the source sheet.
$ws.Range("C2:C13").NumberFormat = "0.00E+00"

# Column E needs to be a bit wider to fit the destination-name text.
$ws.Columns.Item(5).ColumnWidth = 11.83

# --- Tweak the selection on the "energy_equivs" sheet (now pushed one tab
# to the right) to match what was left selected in the saved file --------
$equivs = $wb.Worksheets.Item("energy_equivs")
$equivs.Activate()
$equivs.Range("D4").Select()

# --- Leave focus on the new "time" tab, matching the saved selection ------
$ws.Activate()
$ws.Range("E5").Select()
